# Edit News and FeaturedPubs sheets per commit: "made updates to index etc"

$wb = $excel.ActiveWorkbook

# --- News sheet: refresh the two news rows, newest first ---
$news = $wb.Worksheets.Item("News")

$news.Range("A2").Value = 45992
$news.Range("B2").Value = 'Rish Prakash (SMVDU Katra; Intern at IIT Goa), Dr. Anuj Abraham (Research Fellow, IIT Goa), and Dr. Shitala Prasad (Assistant Professor, IIT Goa) on receiving the Best Application Paper Award at the International Intelligent Computing and Technology Conference (ICTCon) 2025, held on 2–3 December 2025. Their award-winning paper, “Hybrid Multi-view 3D Object Detection from 2D Images: Fusion of Structure-from-Motion and Learned Depth Priors,” showcases impactful research and innovation in the field of intelligent computing and computer vision'
$news.Range("C2").Value = 'https://ictcon2025.cit.ac.in/'

$news.Range("A3").Value = 45962
$news.Range("B3").Value = '[Student Achievement] Vipin Gautam, PhD student in the Department of Computer Science and Engineering, received the Best Oral Presentation Award at the Goa Research Scholars Meet 2025, organized by the Centre for Research, Development & Innovation (RDI), Department of Higher Education, Goa, on 20–21 November 2025. He presented his work titled “Bridging RGB–IR Domains for Aerial Object Detection.” Vipin is supervised by Dr. Shitala Prasad, Assistant Professor, Computer Science and Engineering and Dr. Sharad Sinha, Associate Professor, Computer Science and Engineering.'
$news.Range("C3").Value = '#'

# --- FeaturedPubs sheet: swap in the new featured publication ---
$pubs = $wb.Worksheets.Item("FeaturedPubs")

$pubs.Range("A2").Value = 'IEEE SPL'
$pubs.Range("B2").Value = 'SequenceOut: Boosting CNNs by Freezing Layers'
$pubs.Range("C2").Value = 'S Prasad, R Paul, M Kamat'

# --- column width on FeaturedPubs column B (title column) ---
$pubs.Columns.Item(2).ColumnWidth = 58.85546875

# --- Selections / active sheet: FeaturedPubs is now the one in front ---
$news.Range("C26").Select()
$pubs.Range("F12").Select()
$pubs.Activate()

# --- Page setup on FeaturedPubs (paper size + orientation) ---
$pubs.PageSetup.PaperSize = 9
$pubs.PageSetup.Orientation = 1

Write-Output "done"
